$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10: stock codes refreshed, lookup failed with "Error" for all of them,
# so amount/price columns are cleared out.
$errorCodes = @(
    "900.70.061",
    "900.70.063",
    "900.70.062",
    "900.70.066",
    "900.70.068",
    "900.70.382",
    "900.70.383",
    "900.70.384",
    "900.70.386"
)

$row = 2
foreach ($code in $errorCodes) {
    $ws.Cells.Item($row, 1).Value = $code
    $ws.Cells.Item($row, 2).ClearContents()
    $ws.Cells.Item($row, 3).Value = "Error"
    $ws.Cells.Item($row, 4).ClearContents()
    $ws.Cells.Item($row, 5).ClearContents()
    $ws.Cells.Item($row, 6).ClearContents()
    $row++
}

# Row 11: stock code refreshed too, but the lookup succeeded this time with new values.
$ws.Cells.Item(11, 1).Value = "900.70.388"
$ws.Cells.Item(11, 2).Value = 7
$ws.Cells.Item(11, 3).Value = "stokta mevcut"
$ws.Cells.Item(11, 4).Value = "1.056,73 TL"
$ws.Cells.Item(11, 5).Value = "704,49 TL"
$ws.Cells.Item(11, 6).Value = "915,83 TL"
